$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.273058891296387
$ws.Range("B1").Value = 1.942742228507996
$ws.Range("C1").Value = 2.631101131439209
$ws.Range("D1").Value = 3.749205112457275
$ws.Range("E1").Value = 1.08298647403717
